# Apply the "newest EPS-US files" update to the FoTOMRAEL workbook.
#
# Summary of the change (see commit "Uploading newest EPS-US files"):
#  - About sheet: the old "Sources:" block (US EPS / Version 3.1.3 / blank /
#    European Commission / 2017 / citation / hyperlink / "p. 17") is replaced
#    by a short "none specified" note, and the trailing "US vs EU EPS" note
#    at the bottom of the page is removed. The rest of the descriptive text
#    is unchanged, just shifted up.
#  - FoTOMRAEL sheet: the value of "Share of Technology Outside Modeled
#    Region" changes from 0.2 (20%) to 0.25 (25%).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("About")
$ws2 = $wb.Worksheets.Item("FoTOMRAEL")

# --- About sheet -----------------------------------------------------

# B9 currently holds the hyperlinked citation and carries the "Hyperlink"
# cell style; capture that formatting onto B6 (which becomes the new last
# line of the trimmed-down Sources block) before anything gets shifted
# around or the hyperlink is removed.
$ws1.Range("B9").Copy()
$ws1.Range("B6").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws1.Range("B6").ClearContents()

# Remove the old "2017 / citation / hyperlink / p. 17" rows (7-10), then
# re-insert a single blank spacer row so the body text below lands exactly
# one row higher than before (rows 12-30 -> rows 9-27).
$ws1.Rows("7:10").Delete()
$ws1.Rows("8:8").Insert()

# Remove the trailing "US vs EU EPS" paragraph (now at rows 29-31).
$ws1.Rows("29:31").Delete()

# Clear the old "Version 3.1.3" placeholder text.
$ws1.Range("B4").ClearContents()

# Drop the now-orphaned hyperlink object entirely.
$ws1.Hyperlinks.Delete()

# Replace the "US EPS" source label with a simple placeholder note.
$ws1.Range("B3").Value = "none specified"

# --- FoTOMRAEL sheet ---------------------------------------------------

# Share of Technology Outside Modeled Region: 20% -> 25%.
$ws2.Range("B2").Value = 0.25
